$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has: A=MSSV, B=Ho ten, C=Email, D=STT(1,2,3).
# Target layout adds a "Lop" (class) column between B and C, so the old
# Email column (C) shifts to D and the old STT column (D) shifts to E.
# Inserting a blank column right before the existing D (i.e. before the
# Email column's right-hand neighbour) keeps C (Email) untouched in place
# and pushes the old D (STT numbers) out to E - exactly the shift needed,
# while leaving the class column to be filled in at C afterwards.
$ws.Columns("D").Insert()

# Move the e-mail addresses (currently still in C) into the new D column.
$ws.Range("D1").Value2 = $ws.Range("C1").Value2
$ws.Range("D2").Value2 = $ws.Range("C2").Value2
$ws.Range("D3").Value2 = $ws.Range("C3").Value2

# C becomes the class ("Lop") column for every student row.
$ws.Range("C1").Value2 = "21DTHA1"
$ws.Range("C2").Value2 = "21DTHA1"
$ws.Range("C3").Value2 = "21DTHA1"

# Former column D (STT 1,2,3), now shifted to E, becomes a constant 3.
$ws.Range("E1").Value2 = 3
$ws.Range("E2").Value2 = 3
$ws.Range("E3").Value2 = 3

# New column D should render with the same width as column C (Email's
# original best-fit width carries across both columns).
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Match the recorded selection in the edited workbook.
$ws.Range("C3").Select()
